# adding NA to lookup table and adjusting enums accordingly
#
# Inserts two new lookup rows ("Not applicable" = -1, "Not available" = 0)
# right after the header row, pushing the existing 131 data rows down by two
# rows (old row 2 -> row 4, ... old row 132 -> row 134), and tightens column
# B's width a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 132

# Shift the existing data rows (2..132) down to (4..134), working from the
# bottom up so we never overwrite a row before it has been read. Doing this
# via plain value copies (instead of Rows.Insert) avoids Excel silently
# manufacturing an extra, unused cell style from the row-insert format copy.
for ($r = $lastRow; $r -ge 2; $r--) {
    $srcA = $ws.Cells.Item($r, 1).Value2
    $srcB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 2, 1).Value = $srcA
    $ws.Cells.Item($r + 2, 2).Value = $srcB
}

# Row 3: areaId = 0, displayName = "Not available" (no border, like the
# regular data rows).
$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = "Not available"

$ws.Range("A1").Copy()
$ws.Cells.Item(3, 1).PasteSpecial(-4122)
$ws.Cells.Item(3, 1).Font.Bold = $false
$ws.Cells.Item(3, 1).Borders.LineStyle = 0
$ws.Cells.Item(3, 1).HorizontalAlignment = -4152
$ws.Cells.Item(3, 1).VerticalAlignment = -4160

$ws.Range("A1").Copy()
$ws.Cells.Item(3, 2).PasteSpecial(-4122)
$ws.Cells.Item(3, 2).Font.Bold = $false
$ws.Cells.Item(3, 2).Borders.LineStyle = 0
$ws.Cells.Item(3, 2).HorizontalAlignment = -4131
$ws.Cells.Item(3, 2).VerticalAlignment = -4160

# Row 2: areaId = -1, displayName = "Not applicable" (keeps the thin border
# that used to sit under the header row).
$ws.Cells.Item(2, 1).Value = -1
$ws.Cells.Item(2, 2).Value = "Not applicable"

$ws.Range("A1").Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4122)
$ws.Cells.Item(2, 1).Font.Bold = $false
$ws.Cells.Item(2, 1).HorizontalAlignment = -4152
$ws.Cells.Item(2, 1).VerticalAlignment = -4160

$ws.Range("A1").Copy()
$ws.Cells.Item(2, 2).PasteSpecial(-4122)
$ws.Cells.Item(2, 2).Font.Bold = $false
$ws.Cells.Item(2, 2).HorizontalAlignment = -4131
$ws.Cells.Item(2, 2).VerticalAlignment = -4160

# Narrow column B slightly and reset the view to the top of the sheet with
# B2 selected.
$ws.Columns("B").ColumnWidth = 35.8
$ws.Range("B2").Select()

Write-Host "done"
